$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (full swap)
$ws.Range("B16").Value = 6441905
$ws.Range("C16").Value = 'Germany Landesliga'
$ws.Range("D16").Value = 'Germany Landesliga'
$ws.Range("F16").Value = 'SV Helpenstein'
$ws.Range("G16").Value = 'FC Germania Teveren'
$ws.Range("H16").Value = 4
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 'H'
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 3.25
$ws.Range("M16").Value = 2.1
$ws.Range("N16").Value = 3
$ws.Range("O16").Value = 3.25
$ws.Range("P16").Value = 2.1
$ws.Range("Q16").Value = 0.25
$ws.Range("R16").Value = 1.9
$ws.Range("S16").Value = 1.9
$ws.Range("T16").Value = 3.25
$ws.Range("U16").Value = 1.9
$ws.Range("V16").Value = 1.9
$ws.Range("W16").Value = 2
$ws.Range("X16").Value = -1
$ws.Range("Y16").Value = -1
$ws.Range("Z16").Value = 0.8999999999999999
$ws.Range("AA16").Value = -1
$ws.Range("AB16").Value = 0.8999999999999999
$ws.Range("AC16").Value = -1

# Row 17 (full swap)
$ws.Range("B17").Value = 6441941
$ws.Range("C17").Value = 'Germany Landesliga'
$ws.Range("D17").Value = 'Germany Landesliga'
$ws.Range("F17").Value = 'BSV Schuren'
$ws.Range("G17").Value = 'Turkspor Dortmund 2000'
$ws.Range("H17").Value = 2
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 'A'
$ws.Range("K17").Value = 5
$ws.Range("L17").Value = 5
$ws.Range("M17").Value = 1.4
$ws.Range("N17").Value = 7.5
$ws.Range("O17").Value = 5.25
$ws.Range("P17").Value = 1.25
$ws.Range("Q17").Value = 1.75
$ws.Range("R17").Value = 1.95
$ws.Range("S17").Value = 1.85
$ws.Range("T17").Value = 3.5
$ws.Range("U17").Value = 2
$ws.Range("V17").Value = 1.8
$ws.Range("W17").Value = -1
$ws.Range("X17").Value = -1
$ws.Range("Y17").Value = 0.25
$ws.Range("Z17").Value = -0.5
$ws.Range("AA17").Value = 0.425
$ws.Range("AB17").Value = 1
$ws.Range("AC17").Value = -1

# Row 49 (full swap)
$ws.Range("B49").Value = 7035046
$ws.Range("C49").Value = 'Germany Landesliga'
$ws.Range("D49").Value = 'Germany Landesliga'
$ws.Range("F49").Value = 'Cronenberger SC'
$ws.Range("G49").Value = 'FC Viersen'
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 2
$ws.Range("J49").Value = 'A'
$ws.Range("K49").Value = 2
$ws.Range("L49").Value = 3.6
$ws.Range("M49").Value = 3
$ws.Range("N49").Value = 2
$ws.Range("O49").Value = 3.6
$ws.Range("P49").Value = 3
$ws.Range("Q49").Value = -0.25
$ws.Range("R49").Value = 1.8
$ws.Range("S49").Value = 2
$ws.Range("T49").Value = 2.75
$ws.Range("U49").Value = 1.8
$ws.Range("V49").Value = 2
$ws.Range("W49").Value = -1
$ws.Range("X49").Value = -1
$ws.Range("Y49").Value = 2
$ws.Range("Z49").Value = -1
$ws.Range("AA49").Value = 1
$ws.Range("AB49").Value = -1
$ws.Range("AC49").Value = 1

# Row 50 (full swap)
$ws.Range("B50").Value = 7035047
$ws.Range("C50").Value = 'Germany Landesliga'
$ws.Range("D50").Value = 'Germany Landesliga'
$ws.Range("F50").Value = 'SC Dsseldorf West'
$ws.Range("G50").Value = 'VfL Viktoria JuchenGarzweiler'
$ws.Range("H50").Value = 3
$ws.Range("I50").Value = 4
$ws.Range("J50").Value = 'A'
$ws.Range("K50").Value = 1.909
$ws.Range("L50").Value = 3.75
$ws.Range("M50").Value = 3.1
$ws.Range("N50").Value = 2.2
$ws.Range("O50").Value = 3.6
$ws.Range("P50").Value = 2.625
$ws.Range("Q50").Value = -0.25
$ws.Range("R50").Value = 2
$ws.Range("S50").Value = 1.8
$ws.Range("T50").Value = 3
$ws.Range("U50").Value = 1.825
$ws.Range("V50").Value = 1.975
$ws.Range("W50").Value = -1
$ws.Range("X50").Value = -1
$ws.Range("Y50").Value = 1.625
$ws.Range("Z50").Value = -1
$ws.Range("AA50").Value = 0.8
$ws.Range("AB50").Value = 0.825
$ws.Range("AC50").Value = -1

# Row 51 (full swap)
$ws.Range("B51").Value = 7089910
$ws.Range("C51").Value = 'Germany Landesliga'
$ws.Range("D51").Value = 'Germany Landesliga'
$ws.Range("F51").Value = 'ASV Mettmann'
$ws.Range("G51").Value = 'TuRU Dsseldorf'
$ws.Range("H51").Value = 2
$ws.Range("I51").Value = 1
$ws.Range("J51").Value = 'H'
$ws.Range("K51").Value = 3.25
$ws.Range("L51").Value = 4
$ws.Range("M51").Value = 1.8
$ws.Range("N51").Value = 3.25
$ws.Range("O51").Value = 4
$ws.Range("P51").Value = 1.8
$ws.Range("Q51").Value = 0.5
$ws.Range("R51").Value = 1.975
$ws.Range("S51").Value = 1.825
$ws.Range("T51").Value = 3.25
$ws.Range("U51").Value = 1.85
$ws.Range("V51").Value = 1.95
$ws.Range("W51").Value = 2.25
$ws.Range("X51").Value = -1
$ws.Range("Y51").Value = -1
$ws.Range("Z51").Value = 0.9750000000000001
$ws.Range("AA51").Value = -1
$ws.Range("AB51").Value = -0.5
$ws.Range("AC51").Value = 0.475

# Row 52 (full swap)
$ws.Range("B52").Value = 7089911
$ws.Range("C52").Value = 'Germany Landesliga'
$ws.Range("D52").Value = 'Germany Landesliga'
$ws.Range("F52").Value = 'Spvgg Steele 0309'
$ws.Range("G52").Value = 'VfB Frohnhausen'
$ws.Range("H52").Value = 4
$ws.Range("I52").Value = 3
$ws.Range("J52").Value = 'H'
$ws.Range("K52").Value = 2.25
$ws.Range("L52").Value = 3.75
$ws.Range("M52").Value = 2.5
$ws.Range("N52").Value = 2.25
$ws.Range("O52").Value = 3.75
$ws.Range("P52").Value = 2.5
$ws.Range("Q52").Value = 0
$ws.Range("R52").Value = 1.8
$ws.Range("S52").Value = 2
$ws.Range("T52").Value = 3.5
$ws.Range("U52").Value = 1.8
$ws.Range("V52").Value = 2
$ws.Range("W52").Value = 1.25
$ws.Range("X52").Value = -1
$ws.Range("Y52").Value = -1
$ws.Range("Z52").Value = 0.8
$ws.Range("AA52").Value = -1
$ws.Range("AB52").Value = 0.8
$ws.Range("AC52").Value = -1

# New row 93
$ws.Range("A92").Copy($ws.Range("A93"))
$ws.Range("E92").Copy($ws.Range("E93"))
$ws.Range("A93").Value = 91
$ws.Range("B93").Value = 7847639
$ws.Range("C93").Value = 'Germany Landesliga'
$ws.Range("D93").Value = 'Germany Landesliga'
$ws.Range("E93").Value2 = 45340.41666666666
$ws.Range("F93").Value = 'VfB Fortuna Chemnitz'
$ws.Range("G93").Value = 'SG Taucha 99'
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 2
$ws.Range("J93").Value = 'A'
$ws.Range("K93").Value = 2.15
$ws.Range("L93").Value = 3.75
$ws.Range("M93").Value = 2.6
$ws.Range("N93").Value = 2.15
$ws.Range("O93").Value = 3.75
$ws.Range("P93").Value = 2.625
$ws.Range("Q93").Value = -0.25
$ws.Range("R93").Value = 1.975
$ws.Range("S93").Value = 1.825
$ws.Range("T93").Value = 2.75
$ws.Range("U93").Value = 1.95
$ws.Range("V93").Value = 1.85
$ws.Range("W93").Value = -1
$ws.Range("X93").Value = -1
$ws.Range("Y93").Value = 1.625
$ws.Range("Z93").Value = -1
$ws.Range("AA93").Value = 0.825
$ws.Range("AB93").Value = -1
$ws.Range("AC93").Value = 0.8500000000000001
